# Apply the "Detect emotion" feature test-case updates to the
# attendance-marking test-plan worksheet.
#
# Summary of the edit:
#   - Five existing "Expected Result" cells get a trailing "(Pass)" marker
#     appended, now that the behaviour has been verified.
#   - Two standalone "(Pass)" markers are added next to TC03 and TC08.
#   - Two brand-new test cases (TC09 "Detect emotion" and TC10 "Detect
#     correct emotion") are appended as new rows at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark existing, now-passing expectations -------------------------------
$ws.Range("D2").Value = "The selected course should be displayed in the dropdown.(Pass)"
$ws.Range("D5").Value = "The session dropdown should populate with valid session values related to the selected course.(Pass)"
$ws.Range("G7").Value = "(Pass)"
$ws.Range("D9").Value = "A validation error message should appear, indicating that duplicate attendance is not allowed for the same session.(Pass)"
$ws.Range("D11").Value = "Attendance should be successfully marked for each student. No duplicates should be allowed.(Pass)"
$ws.Range("D16").Value = " Appropriate error messages should be displayed. The system should prevent invalid actions.(Pass)"
$ws.Range("F18").Value = "(Pass)"

# --- TC09: Detect emotion ---------------------------------------------------
$ws.Range("A19").Value = "TC09"
$ws.Range("B19").Value = "Detect emotion"
$ws.Range("C19").Value = "1. Open attendance page"
$ws.Range("C20").Value = "2. Mark attendance for logged in user."
$ws.Range("D20").Value = "While marking attendance , it should detect emotion of face. (Pass)"

# --- TC10: Detect correct emotion ------------------------------------------
$ws.Range("A21").Value = "TC10"
$ws.Range("B21").Value = "Detect correct emotion"
$ws.Range("C21").Value = "1. Mark attendance for logged in user."
$ws.Range("D21").Value = "Emotion detection can mark incorrect emotion based on Webcam quality and different operation system (can Fail)"

# Leave the selection where the author left it when saving.
$ws.Range("D21").Select() | Out-Null
